$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update revised AgTests (F) / AgPosit (G) values for existing rows ---
$ws.Range("F306").Value = 76607
$ws.Range("F307").Value = 75687
$ws.Range("G307").Value = 6331
$ws.Range("F309").Value = 77845
$ws.Range("G309").Value = 5527
$ws.Range("F310").Value = 79527
$ws.Range("F312").Value = 28529
$ws.Range("F313").Value = 76654
$ws.Range("G313").Value = 3460
$ws.Range("F314").Value = 65403
$ws.Range("G314").Value = 3143
$ws.Range("F315").Value = 56832
$ws.Range("G315").Value = 2628
$ws.Range("F325").Value = 774560
$ws.Range("F332").Value = 485190
$ws.Range("F337").Value = 105617
$ws.Range("F341").Value = 283756
$ws.Range("F344").Value = 136256
$ws.Range("G344").Value = 2474
$ws.Range("F345").Value = 292963
$ws.Range("F350").Value = 128025
$ws.Range("F352").Value = 307152
$ws.Range("F356").Value = 161023
$ws.Range("F358").Value = 159461
$ws.Range("F361").Value = 333276
$ws.Range("F363").Value = 189274
$ws.Range("F365").Value = 185031
$ws.Range("F376").Value = 222310
$ws.Range("F377").Value = 176613
$ws.Range("F379").Value = 180245
$ws.Range("F383").Value = 221315
$ws.Range("F384").Value = 172108
$ws.Range("F395").Value = 752302
$ws.Range("F398").Value = 298795
$ws.Range("F405").Value = 174424
$ws.Range("F423").Value = 439580
$ws.Range("G423").Value = 636
$ws.Range("F430").Value = 168996
$ws.Range("F432").Value = 116904
$ws.Range("G432").Value = 410
$ws.Range("F454").Value = 50177
$ws.Range("F456").Value = 47663
$ws.Range("F457").Value = 75324
$ws.Range("G457").Value = 124
$ws.Range("F461").Value = 43472
$ws.Range("F462").Value = 41940
$ws.Range("F464").Value = 69785
$ws.Range("F465").Value = 57812
$ws.Range("F466").Value = 49130
$ws.Range("F467").Value = 49638
$ws.Range("F468").Value = 39463
$ws.Range("G468").Value = 43

# --- Append new row 469 with the latest daily stats ---
$ws.Range("A469").Value = 44363
$ws.Range("A469").NumberFormat = "yyyy-mm-dd"
$ws.Range("B469").Value = 391210
$ws.Range("C469").Value = 4372
$ws.Range("D469").Value = 61
$ws.Range("E469").Value = 12464
$ws.Range("F469").Value = 30958
$ws.Range("G469").Value = 63
